# Packing list: vertically merge the repeated Item/Description columns
# (B, C, D) for each of the three shipment blocks so the same value is
# shown once per block instead of being repeated on every row.
#
# Blocks (top row holds the value, bottom row is the last data row before
# the "SUB TOTAL:" row):
#   Block 1: rows 22-34
#   Block 2: rows 38-50
#   Block 3: rows 54-65

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packing list")

$blocks = @(
    @{ First = 22; Last = 34 },
    @{ First = 38; Last = 50 },
    @{ First = 54; Last = 65 }
)

foreach ($block in $blocks) {
    $first = $block.First
    $last = $block.Last

    foreach ($col in @("B", "C", "D")) {
        $rng = $ws.Range("$col$first`:$col$last")
        $rng.Merge() | Out-Null
    }
}
